# Natmi LR-pairs update (Sema3a-Plxna4): add the "ECs" sending/target
# cluster to the FAPs/sCs cross so rows 2-10 cover every combination of
# (ECs, FAPs, sCs) x (ECs, FAPs, sCs) for the Sema3a/Plxna4 edge, and
# refresh all of the expression/specificity statistics accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Sema3a", "Plxna4", "ECs", 2, 0.6666666666666666, 1.397441, 4.192323, 0.6676161521996591, 0.6676161521996592, 3, 1, 2.172908, 6.518724, 0.6836640069168457, 0.6836640069168457, 3.036510728428, 27.328596555852, 0.4564251336952256, 0.4564251336952257),
    @("ECs", "Sema3a", "Plxna4", "FAPs", 2, 0.6666666666666666, 1.397441, 4.192323, 0.6676161521996591, 0.6676161521996592, 3, 1, 0.6730553333333334, 2.019166, 0.2117640075251322, 0.2117640075251322, 0.9405551180686668, 8.464996062618, 0.1413770718783084, 0.1413770718783084),
    @("ECs", "Sema3a", "Plxna4", "sCs", 2, 0.6666666666666666, 1.397441, 4.192323, 0.6676161521996591, 0.6676161521996592, 3, 1, 0.332364, 0.9970920000000001, 0.104571985558022, 0.104571985558022, 0.464459080524, 4.180131724716, 0.06981394662612499, 0.069813946626125),
    @("FAPs", "Sema3a", "Plxna4", "ECs", 3, 1, 0.2347746666666667, 0.7043240000000001, 0.112161700990566, 0.112161700990566, 3, 1, 2.172908, 6.518724, 0.6836640069168457, 0.6836640069168457, 0.5101437513973334, 4.591293762576, 0.07668091792181951, 0.07668091792181951),
    @("FAPs", "Sema3a", "Plxna4", "FAPs", 3, 1, 0.2347746666666667, 0.7043240000000001, 0.112161700990566, 0.112161700990566, 3, 1, 0.6730553333333334, 2.019166, 0.2117640075251322, 0.2117640075251322, 0.1580163415315556, 1.422147073784, 0.02375181129259785, 0.02375181129259785),
    @("FAPs", "Sema3a", "Plxna4", "sCs", 3, 1, 0.2347746666666667, 0.7043240000000001, 0.112161700990566, 0.112161700990566, 3, 1, 0.332364, 0.9970920000000001, 0.104571985558022, 0.104571985558022, 0.07803064731200002, 0.7022758258080001, 0.01172897177614866, 0.01172897177614866),
    @("sCs", "Sema3a", "Plxna4", "ECs", 3, 1, 0.4609646666666667, 1.382894, 0.2202221468097748, 0.2202221468097748, 3, 1, 2.172908, 6.518724, 0.6836640069168457, 0.6836640069168457, 1.001633811917333, 9.014704307256, 0.1505579552998005, 0.1505579552998005),
    @("sCs", "Sema3a", "Plxna4", "FAPs", 3, 1, 0.4609646666666667, 1.382894, 0.2202221468097748, 0.2202221468097748, 3, 1, 0.6730553333333334, 2.019166, 0.2117640075251322, 0.2117640075251322, 0.3102547273782223, 2.792292546404, 0.04663512435422591, 0.04663512435422592),
    @("sCs", "Sema3a", "Plxna4", "sCs", 3, 1, 0.4609646666666667, 1.382894, 0.2202221468097748, 0.2202221468097748, 3, 1, 0.332364, 0.9970920000000001, 0.104571985558022, 0.104571985558022, 0.153208060472, 1.378872544248, 0.02302906715574838, 0.02302906715574838)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}
